$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 5).Value = 1
$ws.Cells.Item(2, 6).Value = 0.3333333333333333
$ws.Cells.Item(2, 7).Value = 0.3399353333333333
$ws.Cells.Item(2, 8).Value = 1.019806
$ws.Cells.Item(2, 9).Value = 0.09929991924017606
$ws.Cells.Item(2, 10).Value = 0.09929991924017606
$ws.Cells.Item(2, 13).Value = 61.04160633333334
$ws.Cells.Item(2, 14).Value = 183.124819
$ws.Cells.Item(2, 15).Value = 0.2043613460574534
$ws.Cells.Item(2, 16).Value = 0.2043613460574534
$ws.Cells.Item(2, 17).Value = 20.75019879612378
$ws.Cells.Item(2, 18).Value = 186.751789165114
$ws.Cells.Item(2, 19).Value = 0.0202930651593188
$ws.Cells.Item(2, 20).Value = 0.0202930651593188
$ws.Cells.Item(3, 5).Value = 1
$ws.Cells.Item(3, 6).Value = 0.3333333333333333
$ws.Cells.Item(3, 7).Value = 0.3399353333333333
$ws.Cells.Item(3, 8).Value = 1.019806
$ws.Cells.Item(3, 9).Value = 0.09929991924017606
$ws.Cells.Item(3, 10).Value = 0.09929991924017606
$ws.Cells.Item(3, 15).Value = 0.3559304658284363
$ws.Cells.Item(3, 16).Value = 0.3559304658284363
$ws.Cells.Item(3, 17).Value = 36.14004343786533
$ws.Cells.Item(3, 18).Value = 325.260390940788
$ws.Cells.Item(3, 19).Value = 0.03534386651188197
$ws.Cells.Item(3, 20).Value = 0.03534386651188198
$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 0.3333333333333333
$ws.Cells.Item(4, 7).Value = 0.3399353333333333
$ws.Cells.Item(4, 8).Value = 1.019806
$ws.Cells.Item(4, 9).Value = 0.09929991924017606
$ws.Cells.Item(4, 10).Value = 0.09929991924017606
$ws.Cells.Item(4, 13).Value = 131.3384093333333
$ws.Cells.Item(4, 14).Value = 394.015228
$ws.Cells.Item(4, 15).Value = 0.4397081881141102
$ws.Cells.Item(4, 16).Value = 0.4397081881141103
$ws.Cells.Item(4, 17).Value = 44.64656595619643
$ws.Cells.Item(4, 18).Value = 401.819093605768
$ws.Cells.Item(4, 19).Value = 0.04366298756897529
$ws.Cells.Item(4, 20).Value = 0.04366298756897529
$ws.Cells.Item(5, 9).Value = 0.4094685684206303
$ws.Cells.Item(5, 10).Value = 0.4094685684206303
$ws.Cells.Item(5, 13).Value = 61.04160633333334
$ws.Cells.Item(5, 14).Value = 183.124819
$ws.Cells.Item(5, 15).Value = 0.2043613460574534
$ws.Cells.Item(5, 16).Value = 0.2043613460574534
$ws.Cells.Item(5, 17).Value = 85.56456299769722
$ws.Cells.Item(5, 18).Value = 770.0810669792751
$ws.Cells.Item(5, 19).Value = 0.08367954781065848
$ws.Cells.Item(5, 20).Value = 0.08367954781065848
$ws.Cells.Item(6, 9).Value = 0.4094685684206303
$ws.Cells.Item(6, 10).Value = 0.4094685684206303
$ws.Cells.Item(6, 15).Value = 0.3559304658284363
$ws.Cells.Item(6, 16).Value = 0.3559304658284363
$ws.Cells.Item(6, 19).Value = 0.1457423383000579
$ws.Cells.Item(6, 20).Value = 0.1457423383000579
$ws.Cells.Item(7, 9).Value = 0.4094685684206303
$ws.Cells.Item(7, 10).Value = 0.4094685684206303
$ws.Cells.Item(7, 13).Value = 131.3384093333333
$ws.Cells.Item(7, 14).Value = 394.015228
$ws.Cells.Item(7, 15).Value = 0.4397081881141102
$ws.Cells.Item(7, 16).Value = 0.4397081881141103
$ws.Cells.Item(7, 17).Value = 184.1025207962555
$ws.Cells.Item(7, 18).Value = 1656.9226871663
$ws.Cells.Item(7, 19).Value = 0.1800466823099139
$ws.Cells.Item(7, 20).Value = 0.180046682309914
$ws.Cells.Item(8, 7).Value = 1.681642333333333
$ws.Cells.Item(8, 8).Value = 5.044927
$ws.Cells.Item(8, 9).Value = 0.4912315123391937
$ws.Cells.Item(8, 10).Value = 0.4912315123391937
$ws.Cells.Item(8, 13).Value = 61.04160633333334
$ws.Cells.Item(8, 14).Value = 183.124819
$ws.Cells.Item(8, 15).Value = 0.2043613460574534
$ws.Cells.Item(8, 16).Value = 0.2043613460574534
$ws.Cells.Item(8, 17).Value = 102.6501493048015
$ws.Cells.Item(8, 18).Value = 923.8513437432132
$ws.Cells.Item(8, 19).Value = 0.1003887330874762
$ws.Cells.Item(8, 20).Value = 0.1003887330874762
$ws.Cells.Item(9, 7).Value = 1.681642333333333
$ws.Cells.Item(9, 8).Value = 5.044927
$ws.Cells.Item(9, 9).Value = 0.4912315123391937
$ws.Cells.Item(9, 10).Value = 0.4912315123391937
$ws.Cells.Item(9, 15).Value = 0.3559304658284363
$ws.Cells.Item(9, 16).Value = 0.3559304658284363
$ws.Cells.Item(9, 17).Value = 178.7829066713273
$ws.Cells.Item(9, 18).Value = 1609.046160041946
$ws.Cells.Item(9, 19).Value = 0.1748442610164964
$ws.Cells.Item(9, 20).Value = 0.1748442610164965
$ws.Cells.Item(10, 7).Value = 1.681642333333333
$ws.Cells.Item(10, 8).Value = 5.044927
$ws.Cells.Item(10, 9).Value = 0.4912315123391937
$ws.Cells.Item(10, 10).Value = 0.4912315123391937
$ws.Cells.Item(10, 13).Value = 131.3384093333333
$ws.Cells.Item(10, 14).Value = 394.015228
$ws.Cells.Item(10, 15).Value = 0.4397081881141102
$ws.Cells.Item(10, 16).Value = 0.4397081881141103
$ws.Cells.Item(10, 17).Value = 220.8642291275951
$ws.Cells.Item(10, 18).Value = 1987.778062148356
$ws.Cells.Item(10, 19).Value = 0.215998518235221
$ws.Cells.Item(10, 20).Value = 0.215998518235221
